# Apply the "OpenAccountTest" sheet changes described by the commit.
$wb = $excel.ActiveWorkbook

# --- Sheet "AddCustomerTest" (sheet1): just move the selection, it is no
#     longer the active/selected tab once we activate OpenAccountTest below.
$ws1 = $wb.Worksheets.Item("AddCustomerTest")
$ws1.Range("B2").Select()

# --- Sheet "Sheet2" -> rename to "OpenAccountTest" and populate with data.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "OpenAccountTest"

# Row 1: headers
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"

# Row 2: values (write B2 before A2 so the shared-string table picks up
# "Rupee" ahead of "Mukesh Ambani", matching the target ordering)
$ws2.Range("B2").Value = "Rupee"
$ws2.Range("A2").Value = "Mukesh Ambani"

# Auto-fit column A to the new content (mirrors the bestFit/customWidth
# column sizing that appears on the sheet after adding the data).
$ws2.Columns("A:A").AutoFit()

# Select B12 on this sheet and make it the active sheet/tab.
$ws2.Range("B12").Select()
$ws2.Activate()
